$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: add the Drools "Import" statement row
#   B2 = "Import" (keeps its existing style)
#   C2 = "com.myspace.eotworkflow.*" (loses its explicit style -> reverts to Normal)
$ws.Range("B2").Value = "Import"
$ws.Range("C2").Value = "com.myspace.eotworkflow.*"
$ws.Range("C2").Style = "Normal"

# Row 7 (F7): change the ACTION snippet text from the old condition-ish
# placeholder to the real action statement
$ws.Range("F7").Value = "document.setAutogen(`$param)"

# Move/restore the active selection to F6 (matches the saved cursor position)
[void]$ws.Range("F6").Select()

# Widen column F so the longer action snippet fits
$ws.Columns.Item(6).ColumnWidth = 28.33
